$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-3: name / account number
$ws.Range("C2").Value = "Hartmut"
# B3 holds a 16-digit card number; pre-format as Text so the long digit
# string is stored verbatim instead of being coerced into a Number (which
# would lose the literal formatting / risk precision loss on round numbers).
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Row 5: opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 25.10.2024"

# Row 6
$ws.Range("B6").Value = "27.10."
$ws.Range("C6").Value = "28.10."
$ws.Range("D6").Value = "KARTENZ./27.10 REWE RO"
$ws.Range("E6").Value = "100,81-"

# Row 7
$ws.Range("B7").Value = "30.10."
$ws.Range("C7").Value = "31.10."
$ws.Range("E7").Value = "25,20-"

# Row 8
$ws.Range("B8").Value = "03.11."
$ws.Range("C8").Value = "04.11."
$ws.Range("D8").Value = "BURGER KING Schwandorf"
$ws.Range("E8").Value = "37,34-"

# Row 9
$ws.Range("B9").Value = "06.11."
$ws.Range("C9").Value = "07.11."
$ws.Range("D9").Value = "KARTENZAHLUNG JET TANKSTELLE"
$ws.Range("E9").Value = "83,39-"

# Row 10
$ws.Range("B10").Value = "10.11."
$ws.Range("C10").Value = "11.11."
$ws.Range("D10").Value = "RECHNUNG VODAFONE GMBH 57145650"
$ws.Range("E10").Value = "40,26-"

# Row 11: clear the entry (becomes an empty row in the statement)
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Cells.Item(11, 5).VerticalAlignment = -4108
$ws.Cells.Item(11, 5).WrapText = $true

# Row 12: closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 14.11.2024"
$ws.Range("E12").Value = "287,00-"

# Row 13: next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 21.11.2024"
